# Apply updated crypto price/volume data as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.737.60"
$ws.Range("E2").Value = "  +2.31%  "
$ws.Range("D3").Value = "1.873.57"
$ws.Range("E3").Value = "  +2.10%  "
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").Value = "'324.61"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("D7").Value = "'0.4588"
$ws.Range("E7").Value = "  -0.61%  "
$ws.Range("D8").Value = "'0.3851"
$ws.Range("E8").Value = "  -0.30%  "
$ws.Range("D9").Value = "'0.07861"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").Value = "'0.9927"
$ws.Range("D11").Value = "'21.71"
$ws.Range("E11").Value = "  -1.41%  "
$ws.Range("D12").Value = "1.886.39"
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("D13").Value = "'6.965"
$ws.Range("E13").Value = "  +0.96%  "
$ws.Range("D14").Value = "'5.687"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "'0.06970"
$ws.Range("E15").Value = "  +1.84%  "
$ws.Range("D16").Value = "'88.33"
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("E17").Value = "  +0.55%  "
$ws.Range("D18").Value = "'0.00001004"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("D20").Value = "'1.005"
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").Value = "28.746.25"
$ws.Range("E21").Value = "  +2.32%  "
$ws.Range("D22").Value = "'5.270"
$ws.Range("E22").Value = "  -0.40%  "
$ws.Range("D23").Value = "'11.01"
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("D24").Value = "'2.124"
$ws.Range("E24").Value = "  +1.80%  "
$ws.Range("D25").Value = "2.111.07"
$ws.Range("E25").Value = "  +1.77%  "
$ws.Range("D26").Value = "'153.15"
$ws.Range("E26").Value = "  -0.98%  "
$ws.Range("D27").Value = "'19.22"
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("D28").Value = "'5.766"
$ws.Range("E28").Value = "  +0.58%  "
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "'118.85"
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("B30").Value = "LidoDAOToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D30").Value = "'1.940"
$ws.Range("E30").Value = "  -1.90%  "
$ws.Range("D31").Value = "'0.09302"
$ws.Range("E31").Value = "  +0.52%  "
$ws.Range("D32").Value = "'0.9145"
$ws.Range("E32").Value = "  -2.79%  "
$ws.Range("D33").Value = "'5.294"
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("D34").Value = "'1.334"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("D35").Value = "'3.309"
$ws.Range("E35").Value = "  -0.62%  "
$ws.Range("D36").Value = "'0.05750"
$ws.Range("E36").Value = "  -1.38%  "
$ws.Range("D37").Value = "'1.146"
$ws.Range("E37").Value = "  +0.54%  "
$ws.Range("D38").Value = "'0.02070"
$ws.Range("E38").Value = "  -2.04%  "
$ws.Range("D39").Value = "'7.690"
$ws.Range("E39").Value = "  -0.58%  "
$ws.Range("D40").Value = "'0.5620"
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("E41").Value = "  +1.56%  "
$ws.Range("D42").Value = "'9.851"
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("D43").Value = "'0.07183"
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("D44").Value = "'11.76"
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("D45").Value = "'0.5268"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "'2.128"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "'1.116"
$ws.Range("E47").Value = "  -1.27%  "
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("D49").Value = "'113.31"
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("D50").Value = "'2.412"
$ws.Range("E50").Value = "  +4.16%  "
$ws.Range("D51").Value = "'1.005"
$ws.Range("E51").Value = "  +0.45%  "
